# Add data for 2022-10-12: update 2022 (column I) running totals
# across the Citywide Totals, By Neighborhood, and per-neighborhood sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 9).Value = 5696
$ws.Cells.Item(3, 9).Value = 5948
$ws.Cells.Item(4, 9).Value = 1354
$ws.Cells.Item(6, 9).Value = 6674
$ws.Cells.Item(7, 9).Value = 20216

$ws = $wb.Worksheets.Item('Uptown')
$ws.Cells.Item(4, 9).Value = 32
$ws.Cells.Item(6, 9).Value = 85
$ws.Cells.Item(7, 9).Value = 236

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(4, 9).Value = 16
$ws.Cells.Item(7, 9).Value = 218

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Cells.Item(6, 9).Value = 41
$ws.Cells.Item(7, 9).Value = 115

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Cells.Item(6, 9).Value = 20
$ws.Cells.Item(7, 9).Value = 70

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(3, 9).Value = 215
$ws.Cells.Item(4, 9).Value = 36
$ws.Cells.Item(6, 9).Value = 185
$ws.Cells.Item(7, 9).Value = 648

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Cells.Item(3, 9).Value = 137
$ws.Cells.Item(7, 9).Value = 371

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(2, 9).Value = 182
$ws.Cells.Item(3, 9).Value = 290
$ws.Cells.Item(7, 9).Value = 792

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Cells.Item(2, 9).Value = 59
$ws.Cells.Item(6, 9).Value = 78
$ws.Cells.Item(7, 9).Value = 199

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(2, 9).Value = 160
$ws.Cells.Item(4, 9).Value = 79
$ws.Cells.Item(7, 9).Value = 635
$ws.Cells.Item(8, 9).Value = 1215
$ws.Cells.Item(11, 9).Value = 299
$ws.Cells.Item(14, 9).Value = 115
$ws.Cells.Item(15, 9).Value = 228
$ws.Cells.Item(19, 9).Value = 556
$ws.Cells.Item(20, 9).Value = 491
$ws.Cells.Item(23, 9).Value = 202
$ws.Cells.Item(25, 9).Value = 102
$ws.Cells.Item(29, 9).Value = 1270
$ws.Cells.Item(30, 9).Value = 70
$ws.Cells.Item(31, 9).Value = 199
$ws.Cells.Item(33, 9).Value = 919
$ws.Cells.Item(36, 9).Value = 262
$ws.Cells.Item(37, 9).Value = 648
$ws.Cells.Item(42, 9).Value = 682
$ws.Cells.Item(43, 9).Value = 175
$ws.Cells.Item(44, 9).Value = 148
$ws.Cells.Item(47, 9).Value = 139
$ws.Cells.Item(51, 9).Value = 234
$ws.Cells.Item(52, 9).Value = 444
$ws.Cells.Item(53, 9).Value = 211
$ws.Cells.Item(54, 9).Value = 424
$ws.Cells.Item(55, 9).Value = 224
$ws.Cells.Item(57, 9).Value = 80
$ws.Cells.Item(63, 9).Value = 73
$ws.Cells.Item(67, 9).Value = 792
$ws.Cells.Item(71, 9).Value = 61
$ws.Cells.Item(73, 9).Value = 180
$ws.Cells.Item(77, 9).Value = 130
$ws.Cells.Item(78, 9).Value = 274
$ws.Cells.Item(81, 9).Value = 19
$ws.Cells.Item(83, 9).Value = 429
$ws.Cells.Item(85, 9).Value = 914
$ws.Cells.Item(86, 9).Value = 124
$ws.Cells.Item(89, 9).Value = 236
$ws.Cells.Item(90, 9).Value = 249
$ws.Cells.Item(91, 9).Value = 217
$ws.Cells.Item(93, 9).Value = 114
$ws.Cells.Item(96, 9).Value = 218
$ws.Cells.Item(97, 9).Value = 168
$ws.Cells.Item(99, 9).Value = 371
$ws.Cells.Item(101, 9).Value = 20216

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Cells.Item(2, 9).Value = 147
$ws.Cells.Item(4, 9).Value = 17
$ws.Cells.Item(7, 9).Value = 429

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(3, 9).Value = 345
$ws.Cells.Item(6, 9).Value = 288
$ws.Cells.Item(7, 9).Value = 919

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(2, 9).Value = 93
$ws.Cells.Item(7, 9).Value = 424

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 9).Value = 372
$ws.Cells.Item(3, 9).Value = 437
$ws.Cells.Item(6, 9).Value = 351
$ws.Cells.Item(7, 9).Value = 1270

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(6, 9).Value = 162
$ws.Cells.Item(7, 9).Value = 556

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Cells.Item(6, 9).Value = 44
$ws.Cells.Item(7, 9).Value = 148

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 9).Value = 250
$ws.Cells.Item(3, 9).Value = 353
$ws.Cells.Item(7, 9).Value = 914

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(6, 9).Value = 208
$ws.Cells.Item(7, 9).Value = 682

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Cells.Item(2, 9).Value = 64
$ws.Cells.Item(4, 9).Value = 37
$ws.Cells.Item(7, 9).Value = 274

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Cells.Item(2, 9).Value = 67
$ws.Cells.Item(7, 9).Value = 224

$ws = $wb.Worksheets.Item('Douglas')
$ws.Cells.Item(3, 9).Value = 70
$ws.Cells.Item(7, 9).Value = 202

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Cells.Item(6, 9).Value = 60
$ws.Cells.Item(7, 9).Value = 217

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(6, 9).Value = 166
$ws.Cells.Item(7, 9).Value = 491

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Cells.Item(6, 9).Value = 79
$ws.Cells.Item(7, 9).Value = 262

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Cells.Item(2, 9).Value = 31
$ws.Cells.Item(7, 9).Value = 114

$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(2, 9).Value = 121
$ws.Cells.Item(7, 9).Value = 444

$ws = $wb.Worksheets.Item('East Side')
$ws.Cells.Item(2, 9).Value = 36
$ws.Cells.Item(4, 9).Value = 6
$ws.Cells.Item(7, 9).Value = 102

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Cells.Item(4, 9).Value = 13
$ws.Cells.Item(6, 9).Value = 47
$ws.Cells.Item(7, 9).Value = 139

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(2, 9).Value = 70
$ws.Cells.Item(6, 9).Value = 85
$ws.Cells.Item(7, 9).Value = 228

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Cells.Item(3, 9).Value = 58
$ws.Cells.Item(6, 9).Value = 81
$ws.Cells.Item(7, 9).Value = 299

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Cells.Item(4, 9).Value = 17
$ws.Cells.Item(6, 9).Value = 47
$ws.Cells.Item(7, 9).Value = 180

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Cells.Item(3, 9).Value = 52
$ws.Cells.Item(7, 9).Value = 160

$ws = $wb.Worksheets.Item('West Town')
$ws.Cells.Item(2, 9).Value = 28
$ws.Cells.Item(6, 9).Value = 105
$ws.Cells.Item(7, 9).Value = 168

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(3, 9).Value = 341
$ws.Cells.Item(6, 9).Value = 394
$ws.Cells.Item(7, 9).Value = 1215

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Cells.Item(4, 9).Value = 59
$ws.Cells.Item(7, 9).Value = 124

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Cells.Item(2, 9).Value = 82
$ws.Cells.Item(7, 9).Value = 249

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Cells.Item(6, 9).Value = 96
$ws.Cells.Item(7, 9).Value = 234

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Cells.Item(2, 9).Value = 29
$ws.Cells.Item(7, 9).Value = 80

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Cells.Item(6, 9).Value = 98
$ws.Cells.Item(7, 9).Value = 175

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Cells.Item(6, 9).Value = 96
$ws.Cells.Item(7, 9).Value = 211

$ws = $wb.Worksheets.Item('Oakland')
$ws.Cells.Item(6, 9).Value = 16
$ws.Cells.Item(7, 9).Value = 61

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Cells.Item(3, 9).Value = 45
$ws.Cells.Item(7, 9).Value = 130

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(3, 9).Value = 197
$ws.Cells.Item(7, 9).Value = 635

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Cells.Item(6, 9).Value = 23
$ws.Cells.Item(7, 9).Value = 79

$ws = $wb.Worksheets.Item('Sauganash,Forest Glen')
$ws.Cells.Item(3, 9).Value = 6
$ws.Cells.Item(6, 9).Value = 19
